$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that moved forward by one day
# (45188 -> 45189, i.e. 2023-09-19 -> 2023-09-20) for every data row (2..233).
for ($r = 2; $r -le 233; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
